$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 138-143: P column increases by 2 or 3, AA mirrors P, AB mirrors -P
$updates = @(
    @{ Row = 138; P = 2374 },
    @{ Row = 139; P = 2327 },
    @{ Row = 140; P = 2279 },
    @{ Row = 141; P = 2200 },
    @{ Row = 142; P = 2137 },
    @{ Row = 143; P = 2074 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $p = $u.P
    $ws.Cells.Item($r, 16).Value = $p     # column P
    $ws.Cells.Item($r, 27).Value = $p     # column AA
    $ws.Cells.Item($r, 28).Value = -$p    # column AB
}

# Append new row 144 for period 01-07-2021
$newRow = 144

# Force the period label to be stored as text, not auto-converted to a date serial.
$ws.Cells.Item($newRow, 1).Value = "'01-07-2021"
$ws.Cells.Item($newRow, 1).Style = "Normal"

for ($col = 2; $col -le 15; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 0
}
$ws.Cells.Item($newRow, 16).Value = 2007   # column P
for ($col = 17; $col -le 26; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 0
}
$ws.Cells.Item($newRow, 27).Value = 2007   # column AA
$ws.Cells.Item($newRow, 28).Value = -2007  # column AB
